# Apply edits to Junction_Flooding_378 worksheet:
# 1) Update data rows 2-5 with new values
# 2) Delete row 6 (table now has 4 data rows instead of 5)
# 3) Widen several columns from 7 to 8 characters

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values for rows 2 through 5 ---
# Row 2
$ws.Cells.Item(2, 1).Value = 45059.50694444445
$ws.Cells.Item(2, 2).Value = 13.283
$ws.Cells.Item(2, 3).Value = 8.791
$ws.Cells.Item(2, 4).Value = 3.507
$ws.Cells.Item(2, 5).Value = 28.828
$ws.Cells.Item(2, 6).Value = 21.492
$ws.Cells.Item(2, 7).Value = 10.237
$ws.Cells.Item(2, 8).Value = 30.575
$ws.Cells.Item(2, 9).Value = 16.288
$ws.Cells.Item(2, 10).Value = 6.473
$ws.Cells.Item(2, 11).Value = 9.558999999999999
$ws.Cells.Item(2, 12).Value = 11.331
$ws.Cells.Item(2, 13).Value = 12.141
$ws.Cells.Item(2, 14).Value = 3.376
$ws.Cells.Item(2, 15).Value = 10.527
$ws.Cells.Item(2, 16).Value = 14.419
$ws.Cells.Item(2, 17).Value = 9.413
$ws.Cells.Item(2, 18).Value = 2.798
$ws.Cells.Item(2, 19).Value = 1.636
$ws.Cells.Item(2, 20).Value = 152.85
$ws.Cells.Item(2, 21).Value = 29.125
$ws.Cells.Item(2, 22).Value = 9.717000000000001
$ws.Cells.Item(2, 23).Value = 18.774
$ws.Cells.Item(2, 24).Value = 9.568
$ws.Cells.Item(2, 25).Value = 2.805
$ws.Cells.Item(2, 26).Value = 16.47
$ws.Cells.Item(2, 27).Value = 8.583
$ws.Cells.Item(2, 28).Value = 7.902
$ws.Cells.Item(2, 29).Value = 9.369
$ws.Cells.Item(2, 30).Value = 11.722
$ws.Cells.Item(2, 31).Value = 3.066
$ws.Cells.Item(2, 32).Value = 27.695
$ws.Cells.Item(2, 33).Value = 5.181
$ws.Cells.Item(2, 34).Value = 12.147

# Row 3
$ws.Cells.Item(3, 1).Value = 45059.51388888889
$ws.Cells.Item(3, 2).Value = 18.122
$ws.Cells.Item(3, 3).Value = 13.148
$ws.Cells.Item(3, 4).Value = 1.818
$ws.Cells.Item(3, 5).Value = 39.669
$ws.Cells.Item(3, 6).Value = 31.569
$ws.Cells.Item(3, 7).Value = 14.138
$ws.Cells.Item(3, 8).Value = 53.278
$ws.Cells.Item(3, 9).Value = 22.106
$ws.Cells.Item(3, 10).Value = 9.571999999999999
$ws.Cells.Item(3, 11).Value = 14.024
$ws.Cells.Item(3, 12).Value = 15.842
$ws.Cells.Item(3, 13).Value = 16.896
$ws.Cells.Item(3, 14).Value = 4.589
$ws.Cells.Item(3, 15).Value = 14.287
$ws.Cells.Item(3, 16).Value = 20.132
$ws.Cells.Item(3, 17).Value = 12.343
$ws.Cells.Item(3, 18).Value = 1.405
$ws.Cells.Item(3, 19).Value = 1.07
$ws.Cells.Item(3, 20).Value = 210.135
$ws.Cells.Item(3, 21).Value = 39.933
$ws.Cells.Item(3, 22).Value = 13.187
$ws.Cells.Item(3, 23).Value = 26.506
$ws.Cells.Item(3, 24).Value = 13.806
$ws.Cells.Item(3, 25).Value = 2.482
$ws.Cells.Item(3, 26).Value = 26.556
$ws.Cells.Item(3, 27).Value = 11.648
$ws.Cells.Item(3, 28).Value = 10.476
$ws.Cells.Item(3, 29).Value = 12.338
$ws.Cells.Item(3, 30).Value = 16.529
$ws.Cells.Item(3, 31).Value = 1.196
$ws.Cells.Item(3, 32).Value = 48.831
$ws.Cells.Item(3, 33).Value = 7.288
$ws.Cells.Item(3, 34).Value = 16.487

# Row 4
$ws.Cells.Item(4, 1).Value = 45059.52083333334
$ws.Cells.Item(4, 2).Value = 6.134
$ws.Cells.Item(4, 3).Value = 4.345
$ws.Cells.Item(4, 4).Value = 0.96
$ws.Cells.Item(4, 5).Value = 13.58
$ws.Cells.Item(4, 6).Value = 10.386
$ws.Cells.Item(4, 7).Value = 4.746
$ws.Cells.Item(4, 8).Value = 23.322
$ws.Cells.Item(4, 9).Value = 7.562
$ws.Cells.Item(4, 10).Value = 3.205
$ws.Cells.Item(4, 11).Value = 4.5
$ws.Cells.Item(4, 12).Value = 5.413
$ws.Cells.Item(4, 13).Value = 5.842
$ws.Cells.Item(4, 14).Value = 1.576
$ws.Cells.Item(4, 15).Value = 4.888
$ws.Cells.Item(4, 16).Value = 6.845
$ws.Cells.Item(4, 17).Value = 4.422
$ws.Cells.Item(4, 18).Value = 0.882
$ws.Cells.Item(4, 19).Value = 0.511
$ws.Cells.Item(4, 20).Value = 67.102
$ws.Cells.Item(4, 21).Value = 13.887
$ws.Cells.Item(4, 22).Value = 4.511
$ws.Cells.Item(4, 23).Value = 9.051
$ws.Cells.Item(4, 24).Value = 4.637
$ws.Cells.Item(4, 25).Value = 1.044
$ws.Cells.Item(4, 26).Value = 10.946
$ws.Cells.Item(4, 27).Value = 3.985
$ws.Cells.Item(4, 28).Value = 3.686
$ws.Cells.Item(4, 29).Value = 4.328
$ws.Cells.Item(4, 30).Value = 5.614
$ws.Cells.Item(4, 31).Value = 0.745
$ws.Cells.Item(4, 32).Value = 21.747
$ws.Cells.Item(4, 33).Value = 2.4
$ws.Cells.Item(4, 34).Value = 5.642

# Row 5
$ws.Cells.Item(5, 1).Value = 45059.52777777778
$ws.Cells.Item(5, 2).Value = 17.68
$ws.Cells.Item(5, 3).Value = 13.11
$ws.Cells.Item(5, 4).Value = 1.14
$ws.Cells.Item(5, 5).Value = 38.67
$ws.Cells.Item(5, 6).Value = 31.33
$ws.Cells.Item(5, 7).Value = 13.85
$ws.Cells.Item(5, 8).Value = 51.45
$ws.Cells.Item(5, 9).Value = 21.52
$ws.Cells.Item(5, 10).Value = 9.49
$ws.Cells.Item(5, 11).Value = 14
$ws.Cells.Item(5, 12).Value = 15.5
$ws.Cells.Item(5, 13).Value = 16.49
$ws.Cells.Item(5, 14).Value = 4.47
$ws.Cells.Item(5, 15).Value = 13.91
$ws.Cells.Item(5, 16).Value = 19.73
$ws.Cells.Item(5, 17).Value = 11.81
$ws.Cells.Item(5, 18).Value = 0.77
$ws.Cells.Item(5, 19).Value = 0.76
$ws.Cells.Item(5, 20).Value = 204.4
$ws.Cells.Item(5, 21).Value = 38.81
$ws.Cells.Item(5, 22).Value = 12.84
$ws.Cells.Item(5, 23).Value = 26
$ws.Cells.Item(5, 24).Value = 13.64
$ws.Cells.Item(5, 25).Value = 2.13
$ws.Cells.Item(5, 26).Value = 25.37
$ws.Cells.Item(5, 27).Value = 11.34
$ws.Cells.Item(5, 28).Value = 10.09
$ws.Cells.Item(5, 29).Value = 11.87
$ws.Cells.Item(5, 30).Value = 16.24
$ws.Cells.Item(5, 31).Value = 0.54
$ws.Cells.Item(5, 32).Value = 46.66
$ws.Cells.Item(5, 33).Value = 7.19
$ws.Cells.Item(5, 34).Value = 16.05

# --- Remove former row 6 (its data is no longer present) ---
$ws.Rows.Item(6).Delete()

# --- Adjust column widths: several columns grow from 7 to 8 characters wide ---
$ws.Range("B1:C1").ColumnWidth = 7.166666666666667
$ws.Range("G1:G1").ColumnWidth = 7.166666666666667
$ws.Range("K1:M1").ColumnWidth = 7.166666666666667
$ws.Range("O1:O1").ColumnWidth = 7.166666666666667
$ws.Range("Q1:Q1").ColumnWidth = 7.166666666666667
$ws.Range("V1:V1").ColumnWidth = 7.166666666666667
$ws.Range("X1:X1").ColumnWidth = 7.166666666666667
$ws.Range("AA1:AD1").ColumnWidth = 7.166666666666667
$ws.Range("AH1:AH1").ColumnWidth = 7.166666666666667
